$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting rows 52:156 down to 53:157
$ws.Rows("52:52").Insert()

# Populate the new row 52 with the data (same as old row 52 but with updated
# date / price / origin / price-per-kg values)
$ws.Range("A52").Value = 5
$ws.Range("B52").Value = "Macroferia Regional de Talca"
$ws.Range("C52").Value = "Maule"
$ws.Range("D52").Value = 45260
$ws.Range("D52").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = 100112026
$ws.Range("G52").Value = "Haba"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 200
$ws.Range("K52").Value = 13000
$ws.Range("L52").Value = 13000
$ws.Range("M52").Value = 13000
$ws.Range("N52").Value = "`$/saco 25 kilos"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 520
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
